# Update Vic file with substitute past results and non-classic pref flows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Substitute past results: update the "Raw numbers" row for the latest Newspoll (row 5, columns C:G)
$ws.Range("C5").Value = 53
$ws.Range("D5").Value = 56
$ws.Range("E5").Value = 46
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 55

# Non-classic preference flow workings, added below the existing tables
$ws.Range("A38").Formula = "=9.26*0.832"
$ws.Range("A39").Formula = "=20.48*0.619"
$ws.Range("B39").Formula = "=(A39-A38)/(20.48-9.26)"

# Update the active selection to match the author's final cursor position
$ws.Range("L26").Select()
